# Update cryptocurrency price/volume snapshot (+ two re-ranked coins)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number
# by Excel need NumberFormat forced to Text first so the literal string is kept.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.480.62"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.909.91"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "238.89"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4775"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("D8").Value = "0.2828"
$ws.Range("E8").Value = "  -3.45%  "
$ws.Range("D9").Value = "0.06689"
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("D10").Value = "18.59"
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("D11").Value = "100.51"
$ws.Range("E11").Value = "  -4.72%  "
$ws.Range("D12").Value = "1.915.73"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").Value = "0.07686"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "5.181"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "0.6663"
$ws.Range("E15").Value = "  -4.42%  "
$ws.Range("D16").Value = "30.501.90"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "255.77"
$ws.Range("E17").Value = "  -7.08%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "0.000007450"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("D20").Value = "12.61"
$ws.Range("E20").Value = "  -3.78%  "
$ws.Range("D21").Value = "5.360"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "6.271"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").Value = "9.306"
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("D25").Value = "167.23"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").Value = "19.05"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").Value = "2.048"
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "4.768"
$ws.Range("E28").Value = "  +4.54%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.387"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "0.1000"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").Value = "4.254"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").Value = "0.04691"
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("D34").Value = "0.7225"
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").Value = "1.102"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "2.703"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "0.01906"
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("D39").Value = "2.607"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "6.257"
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("D41").Value = "74.61"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").Value = "1.953"
$ws.Range("E42").Value = "  -6.83%  "
$ws.Range("D43").Value = "0.8587"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("D44").Value = "105.45"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "0.4220"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("D47").Value = "7.323"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").Value = "0.1196"
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "910.77"
$ws.Range("E49").Value = "  -8.29%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "34.60"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").Value = "8.723"
$ws.Range("E51").Value = "  -4.75%  "
